$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "01.01.2025"
$ws.Range("B5").Select()
